$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geral")
$ws.Activate()

$ws.Range("B38").Value = 0.0003646722
$ws.Range("C38").Value = 0.0244330261
$ws.Range("D38").Value = 2.0326819075
$ws.Range("E38").Value = 271.69306284
$ws.Range("F38").Value = 29745.6117191
$ws.Range("G38").Value = 739743.55642752
$ws.Range("H38").Value = 298860.776139345

$ws.Range("B39").Value = 0.00072934459
$ws.Range("C39").Value = 0.006928768
$ws.Range("D39").Value = 0.08715661
$ws.Range("E39").Value = 1.101309537
$ws.Range("F39").Value = 13.8852522983
$ws.Range("G39").Value = 72.33087422826
$ws.Range("H39").Value = 154.0616471

$ws.Range("B40").Value = 0.00036467229817
$ws.Range("C40").Value = 0.0083874567
$ws.Range("D40").Value = 0.1068489054
$ws.Range("E40").Value = 1.26541195
$ws.Range("F40").Value = 16.926617046
$ws.Range("G40").Value = 76.8192576006
$ws.Range("H40").Value = 158.270694788249

$ws.Range("B41").Value = 0.00072934459
$ws.Range("C41").Value = 0.00072934459
$ws.Range("D41").Value = 0.61811909459
$ws.Range("E41").Value = 58.7771088071
$ws.Range("F41").Value = 5917.273055477
$ws.Range("G41").Value = 139065.63039
$ws.Range("H41").Value = 558118.936776

$ws.Range("B42").Value = 0.00072934459
$ws.Range("C42").Value = 0.0069287686
$ws.Range("D42").Value = 0.111224969
$ws.Range("E42").Value = 1.789445662043
$ws.Range("F42").Value = 43.34564709228
$ws.Range("G42").Value = 923.6486166166
$ws.Range("H42").Value = 3940.740794766

$ws.Range("B43").Value = 0.00036467229817
$ws.Range("C43").Value = 0.014222209
$ws.Range("D43").Value = 1.012694233
$ws.Range("E43").Value = 96.020334112515
$ws.Range("F43").Value = 9003.37941572246
$ws.Range("G43").Value = 224160.39771131
$ws.Range("H43").Value = 903981.375821059

$ws.Range("E48").Value = 99.283491868537
$ws.Range("F48").Value = 9664.62109818
$ws.Range("G48").Value = 245270.455365966
$ws.Range("H48").Value = 967866.559841994

$ws.Range("B49").Value = 0.00036467229817941
$ws.Range("C49").Value = 0.00036467
$ws.Range("D49").Value = 0.002188033
$ws.Range("E49").Value = 0.0207863
$ws.Range("F49").Value = 0.20749853766
$ws.Range("G49").Value = 1.01707029
$ws.Range("H49").Value = 2.146461147

$ws.Range("B50").Value = 0.000364672
$ws.Range("C50").Value = 0.00765811
$ws.Range("D50").Value = 0.7917035593
$ws.Range("E50").Value = 59.3376652998
$ws.Range("F50").Value = 11124.6650484941
$ws.Range("G50").Value = 140101.38903577
$ws.Range("H50").Value = 1193754.9044777

$ws.Range("B51").Value = 0.000364672
$ws.Range("C51").Value = 0.0040113952
$ws.Range("D51").Value = 0.0266210777
$ws.Range("E51").Value = 0.3650369704
$ws.Range("F51").Value = 4.1280904153
$ws.Range("G51").Value = 24.222610395
$ws.Range("H51").Value = 44.254806745562

$ws.Range("B52").Value = 0.000364672
$ws.Range("C52").Value = 0.000364672
$ws.Range("D52").Value = 0.002917378385
$ws.Range("E52").Value = 0.031361817
$ws.Range("F52").Value = 0.3066894
$ws.Range("G52").Value = 1.527611142
$ws.Range("H52").Value = 3.35206776486

$ws.Range("B53").Value = 0.000364672
$ws.Range("C53").Value = 0.0021880321
$ws.Range("D53").Value = 0.0506894124
$ws.Range("E53").Value = 0.27459804025
$ws.Range("F53").Value = 2.812715384419
$ws.Range("G53").Value = 14.8370463024
$ws.Range("H53").Value = 30.0471520937644

$ws.Range("B54").Value = 0.000364672
$ws.Range("C54").Value = 0.01094016894
$ws.Range("D54").Value = 0.995920046
$ws.Range("E54").Value = 93.853886
$ws.Range("F54").Value = 9437.8098802853
$ws.Range("G54").Value = 226727.231519333
$ws.Range("H54").Value = 949199.76092084

$ws.Range("B59").Value = 0.00072934459
$ws.Range("C59").Value = 0.027350422
$ws.Range("D59").Value = 2.47904228302
$ws.Range("E59").Value = 249.61381203622
$ws.Range("F59").Value = 25154.97916262
$ws.Range("G59").Value = 597291.35909619
$ws.Range("H59").Value = 2347234.46204271

$ws.Range("B60").Value = 0.0003646722
$ws.Range("C60").Value = 0.0361025575
$ws.Range("D60").Value = 3.16243816
$ws.Range("E60").Value = 243.53071343
$ws.Range("F60").Value = 23950.1095475583
$ws.Range("G60").Value = 612276.407178204
$ws.Range("H60").Value = 2402872.13422274

$ws.Range("B61").Value = 0.00072934459
$ws.Range("C61").Value = 0.0145868919
$ws.Range("D61").Value = 1.22931031
$ws.Range("E61").Value = 118.001391589
$ws.Range("F61").Value = 11270.58101049
$ws.Range("G61").Value = 180329.1850674
$ws.Range("H61").Value = 1125942.85658021

$ws.Range("B62").Value = 0.00072934459
$ws.Range("C62").Value = 0.00437606
$ws.Range("D62").Value = 0.02880911155
$ws.Range("E62").Value = 0.34935606
$ws.Range("F62").Value = 4.2743240069
$ws.Range("G62").Value = 23.9582231
$ws.Range("H62").Value = 46.170795

$ws.Range("B63").Value = 0.0003646722
$ws.Range("C63").Value = 0.0131282027
$ws.Range("D63").Value = 1.1811735738
$ws.Range("E63").Value = 115.4023721
$ws.Range("F63").Value = 11766.56925054
$ws.Range("G63").Value = 279936.685641768
$ws.Range("H63").Value = 1180015.54233334

$ws.Range("B64").Value = 0.0003646722
$ws.Range("C64").Value = 0.0003646722
$ws.Range("D64").Value = 0.04011392354
$ws.Range("E64").Value = 0.9649221972219
$ws.Range("F64").Value = 59.8321049963
$ws.Range("G64").Value = 2012.665789022
$ws.Range("H64").Value = 8268.4390943

$ws.Range("B65").Value = 0.00072934459
$ws.Range("C65").Value = 0.012398858138
$ws.Range("D65").Value = 1.0480681849
$ws.Range("E65").Value = 96.008734630885
$ws.Range("F65").Value = 11075.3613537802
$ws.Range("G65").Value = 240086.623829858
$ws.Range("H65").Value = 1004123.06668981

$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H66").Select()
